$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Auditorias de productos" table: insert 3 new rows (Analisis y Diseño,
# Casos de Pruebas, Entrega Proyecto) right after "Plan de proyecto" (row 28),
# before "Presentación y Seguimiento" (old row 29). This pushes every row
# from the old row 29 onward down by 3 (the blank separator row that used to
# sit at row 36 ends up at row 39, etc.).
$ws.Rows("29:31").Insert(-4121, 0)

# The freshly inserted rows don't carry the table's usual formatting, so
# copy it over from the row directly above (row 28) - format only.
$ws.Range("A28:D28").Copy()
$ws.Range("A29:D31").PasteSpecial(-4122)
[void]($excel.CutCopyMode = 0)

# Fill in the 3 new rows.
$ws.Range("A29").Value = "Analisis y Diseño"
$ws.Range("B29").Value = "Al finalizar la ejecución"
$ws.Range("C29").Value = "Ariana Sosa"
$ws.Range("D29").Value = "Analista de requerimientos"

$ws.Range("A30").Value = "Casos de Pruebas"
$ws.Range("B30").Value = "Al finalizar la ejecución"
$ws.Range("C30").Value = "Ariana Sosa"
$ws.Range("D30").Value = "Analista de requerimientos"

$ws.Range("A31").Value = "Entrega Proyecto"
$ws.Range("B31").Value = "Al finalizar la entrega del cliente"
$ws.Range("C31").Value = "Ariana Sosa"
$ws.Range("D31").Value = "Líder de proyecto"

# --- Fix a couple of "Auditado" values that were swapped in the
# "Auditorias de procesos" table.
$ws.Range("D10").Value = "Líder de proyecto"
$ws.Range("D12").Value = "Analista de requerimientos"
$ws.Range("D25").Value = "Líder de proyecto"

# Leave the selection where the author last edited.
[void]$ws.Range("D10").Select()
